$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.477.08'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.575.28'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '291.91'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3736'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.98'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3398'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.147'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07575'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.36'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.017'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.964'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.571.34'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001125'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.96'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06761'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.298'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.35'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.426.71'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.358'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.696'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.13'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.009'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.70'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.747.98'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +5.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.200'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.978'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.838'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.376'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02485'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2289'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06538'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.475'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.34'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6221'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.04'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.812'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5823'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.69'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.077'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.224'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07333'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.02%  '
